# Daily attendance processing - 2025-12-29 08:42:33
# Normalizes the "Recorded By" column (G) so that the literal token
# "System" is always the last entry in the comma-separated list of
# recorders, instead of appearing earlier in the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ","
        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }

        $hasSystem = $false
        foreach ($p in $trimmed) {
            if ($p.Equals("System")) {
                $hasSystem = $true
            }
        }

        if ($hasSystem) {
            $newParts = @()
            foreach ($p in $trimmed) {
                if (-not $p.Equals("System")) {
                    $newParts += $p
                }
            }
            $newParts += "System"
            $newValue = $newParts -join ", "

            if (-not $newValue.Equals($value)) {
                $cell.Value2 = $newValue
            }
        }
    }
}
